$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.024779796600342
$ws.Range("B1").Value = 2.196696519851685
$ws.Range("C1").Value = 2.207910060882568
$ws.Range("D1").Value = 2.748614311218262
$ws.Range("E1").Value = 3.499426126480103
